# Update carjacking monthly YoY data for 2021-09-06 data pull

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (June) - 2021 column (H) value corrected
$ws.Range("H7").Value = 130

# Row 9 (August) - update "through" date label and figures
$ws.Range("A9").Value = "August (through 08-29)"
$ws.Range("B9").Value = 30
$ws.Range("C9").Value = 74
$ws.Range("D9").Value = 84
$ws.Range("E9").Value = 62
$ws.Range("G9").Value = 155
$ws.Range("H9").Value = 148

# Row 10 (Total) - updated totals
$ws.Range("B10").Value = 192
$ws.Range("C10").Value = 376
$ws.Range("D10").Value = 549
$ws.Range("E10").Value = 487
$ws.Range("G10").Value = 776
$ws.Range("H10").Value = 1061
